$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin names and links (column B, C) -- plain text values
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('B10').Value = 'WrappedEther'
$ws.Range('C10').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('B11').Value = 'Solana'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('B21').Value = 'BinanceUSD'
$ws.Range('C21').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('B25').Value = 'BitcoinCash'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('B46').Value = 'Algorand'
$ws.Range('C46').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'

# Price and Volume columns (D, E) -- force text storage so numeric-looking
# strings (e.g. "0.9980", "243.07") keep their exact original formatting
# instead of being auto-parsed into doubles by Excel.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.002.70'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.93%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.880.94'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.72%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.07'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -4.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9980'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4959'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -3.86%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2921'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06638'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.879.61'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.81%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '16.74'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -4.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07237'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6693'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -4.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '86.48'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.879'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.968.02'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007912'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9983'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.79'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.78%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.122.80'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9971'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.772'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.60%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.680'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.066'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.79'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '149.93'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.13'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.918'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.90%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.390'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.61%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.191'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08749'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.944'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05076'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.70%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7131'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.17%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.113'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.665'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.82%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01783'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +5.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.691'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -5.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.177'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -6.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9324'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.53%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.806'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4247'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.68%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9979'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '102.45'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.486'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -3.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1261'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05648'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '32.44'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.273'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3769'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.72%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '55.95'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.60%  '

# Strip the temporary text number-format back off so the cells end up
# with the same (default) style as in the original workbook.
$ws.Range('D2,E2,D3,E3,D5,E5,D6,E6,D7,E7,D8,E8,D9,E9,D10,E10,D11,E11,D12,E12,D13,E13,D14,E14,D15,E15,D16,E16,D17,E17,D18,E18,D19,E19,D20,E20,D21,E21,D22,E22,D23,E23,D24,E24,D25,E25,D26,E26,D27,E27,D28,E28,D29,E29,D30,E30,D31,E31,D32,E32,D33,E33,D34,E34,D35,E35,D36,E36,D37,E37,D38,E38,D39,E39,D40,E40,D41,E41,D42,E42,D43,E43,D44,E44,D45,E45,D46,E46,D47,E47,D48,E48,D49,E49,D50,E50,D51,E51').ClearFormats()
